# Apply cell updates described by the commit diff for cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '51.276.98'
$ws.Cells.Item(2, 5).Value = '  -2.00%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.749.42'
$ws.Cells.Item(3, 5).Value = '  -2.80%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'354.63"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.48%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'107.02"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -4.21%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.548"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -3.48%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'0.999"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -0.07%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.579"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -4.10%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'39.08"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -4.35%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +3.09%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -3.66%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'19.66"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.23%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -4.34%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.180.27'
$ws.Cells.Item(15, 5).Value = '  -3.15%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.747.69'
$ws.Cells.Item(16, 5).Value = '  -3.21%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'0.918"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -1.12%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '51.154.74'
$ws.Cells.Item(18, 5).Value = '  -2.04%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.90%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -4.34%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'12.98"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -2.88%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0960'
$ws.Cells.Item(22, 5).Value = '  -3.87%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'69.33"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -1.40%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'264.94"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -2.80%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'2.74"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -2.52%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.01%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -3.34%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.162"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +13.57%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.69%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -2.49%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = "'6.07"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +3.23%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(32, 4).Value = "'34.49"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.66%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'51.15"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -2.34%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'0.0441"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -7.04%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.0829"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -2.60%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -7.67%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -0.13%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'18.46"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +0.69%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'3.11"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -4.47%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  -4.75%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -3.11%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -2.69%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'2.20"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -3.00%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'119.29"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -4.97%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'21.71"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -3.77%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.077.88'
$ws.Cells.Item(46, 5).Value = '  -0.15%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(47, 4).Value = "'2.31"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -1.01%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'NEARProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(48, 4).Value = "'3.21"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -3.02%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -3.98%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -6.43%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +4.35%  '
